# Add season record columns (Wins, Losses, Ties) to the worksheet,
# using the same header style as the existing header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold, border, centered) onto the
# new header cells so they match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every data row (2-41)
$wins = 93
$losses = 69
$ties = 0

for ($row = 2; $row -le 41; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # AD
    $ws.Cells.Item($row, 31).Value = $losses  # AE
    $ws.Cells.Item($row, 32).Value = $ties    # AF
}
